# Creating new day averages
# Appends 7 new daily tracker rows (rows 291-297, 2023-06-19 through 2023-06-25)
# to the "Main" worksheet of the Dailies workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Data for the new rows. Column A holds the date (serial number), remaining
# columns hold either "X", "XY" or "Y" text markers matching the existing sheet.
$newRows = @(
    @{ Row = 291; Date = 45096; Cells = @{ B = "XY"; C = "X"; D = "X"; F = "X"; G = "X"; H = "Y"; L = "X" } },
    @{ Row = 292; Date = 45097; Cells = @{ B = "XY"; D = "X"; E = "X"; F = "X"; G = "X"; H = "Y"; K = "X"; N = "X" } },
    @{ Row = 293; Date = 45098; Cells = @{ B = "XY"; C = "X"; D = "X"; E = "X"; F = "X"; G = "X"; H = "Y"; I = "X"; L = "X"; N = "X" } },
    @{ Row = 294; Date = 45099; Cells = @{ B = "XY"; D = "X"; F = "X"; G = "X"; H = "Y"; N = "X" } },
    @{ Row = 295; Date = 45100; Cells = @{ B = "XY"; F = "X"; G = "X"; H = "Y"; L = "X"; N = "X" } },
    @{ Row = 296; Date = 45101; Cells = @{ B = "XY"; C = "X"; F = "X"; G = "X"; I = "X"; L = "X"; N = "X" } },
    @{ Row = 297; Date = 45102; Cells = @{ B = "XY"; D = "X"; E = "X"; G = "X"; H = "Y"; K = "X"; L = "X"; N = "X" } }
)

foreach ($rowInfo in $newRows) {
    $r = $rowInfo.Row

    # Date column (A) — numeric serial, formatted the same way as the rows above it
    # (built-in short-date number format, same as the existing date column).
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $rowInfo.Date
    $dateCell.NumberFormat = "m/d/yy"

    foreach ($colLetter in $rowInfo.Cells.Keys) {
        $ws.Range($colLetter + $r).Value = $rowInfo.Cells[$colLetter]
    }
}

# Match the saved selection state from the edit.
$ws.Range("N293").Select()
